$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Overview sheet: Latest HO Xliff Generate Date
$overview.Range("G2").Value = "2016-09-04 21:13:19"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime
$zhcn.Range("H2").Value = "2016-09-04 21:13:14"
$zhcn.Range("K2").Value = "2016-09-04 21:13:31"

# de-de sheet: Correspond Handoff Datetime / Correspond Handback DateTime
$dede.Range("H2").Value = "2016-09-04 21:13:19"
$dede.Range("K2").Value = "2016-09-04 21:13:39"
